# Update crypto price/volume data per the Thu Oct 17 16:48:14 UTC 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.361.22"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "2.614.96"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.04%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.72%  "

$ws.Range("D9").Value = "2.613.33"
$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("E10").Value = "  -2.30%  "

$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").Value = "3.091.06"
$ws.Range("E15").Value = "  +0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.09%  "

$ws.Range("D17").Value = "67.295.03"
$ws.Range("E17").Value = "  -0.43%  "

$ws.Range("D18").Value = "2.619.75"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "367.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.51%  "

$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("E23").Value = "  -1.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "67.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").Value = "2.753.36"
$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "581.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.25%  "

$ws.Range("E34").Value = "  -2.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.125"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.41%  "

$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "155.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.365"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.36%  "

$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("E44").Value = "  +2.46%  "

$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").Value = "0.0₆0294"
$ws.Range("E47").Value = "  +2.55%  "

$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("E49").Value = "  -0.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0787"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.79%  "
